# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E19) previously listed the periods in
# ascending order 2503, 2504, 2505, 2506. The update removes the old
# statement periods and adds the new ones, resulting in the column
# being reordered to 2506, 2505, 2504, 2503.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2506"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2504"
$ws.Range("E19").Value = "2503"
